$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41:E41").Copy()
$ws.Range("A42:E45").PasteSpecial(-4122)

$ws.Cells.Item(12,1).Value = "HRAD"
$ws.Cells.Item(12,2).Value = 2025
$ws.Cells.Item(12,3).Value = 11
$ws.Cells.Item(12,4).Value = 137
$ws.Cells.Item(12,5).Value = 45

$ws.Cells.Item(13,1).Value = "HRAD"
$ws.Cells.Item(13,2).Value = 2025
$ws.Cells.Item(13,3).Value = 1
$ws.Cells.Item(13,4).Value = 103
$ws.Cells.Item(13,5).Value = 60

$ws.Cells.Item(14,1).Value = "HRAD"
$ws.Cells.Item(14,2).Value = 2025
$ws.Cells.Item(14,3).Value = 2
$ws.Cells.Item(14,4).Value = 115
$ws.Cells.Item(14,5).Value = 63

$ws.Cells.Item(15,1).Value = "HRAD"
$ws.Cells.Item(15,2).Value = 2025
$ws.Cells.Item(15,3).Value = 3
$ws.Cells.Item(15,4).Value = 113
$ws.Cells.Item(15,5).Value = 65

$ws.Cells.Item(16,1).Value = "HRAD"
$ws.Cells.Item(16,2).Value = 2025
$ws.Cells.Item(16,3).Value = 4
$ws.Cells.Item(16,4).Value = 118
$ws.Cells.Item(16,5).Value = 70

$ws.Cells.Item(17,1).Value = "HRAD"
$ws.Cells.Item(17,2).Value = 2025
$ws.Cells.Item(17,3).Value = 5
$ws.Cells.Item(17,4).Value = 106
$ws.Cells.Item(17,5).Value = 65

$ws.Cells.Item(18,1).Value = "HRAD"
$ws.Cells.Item(18,2).Value = 2025
$ws.Cells.Item(18,3).Value = 6
$ws.Cells.Item(18,4).Value = 106
$ws.Cells.Item(18,5).Value = 62

$ws.Cells.Item(19,1).Value = "HRAD"
$ws.Cells.Item(19,2).Value = 2025
$ws.Cells.Item(19,3).Value = 7
$ws.Cells.Item(19,4).Value = 109
$ws.Cells.Item(19,5).Value = 64

$ws.Cells.Item(20,1).Value = "HRAD"
$ws.Cells.Item(20,2).Value = 2025
$ws.Cells.Item(20,3).Value = 8
$ws.Cells.Item(20,4).Value = 124
$ws.Cells.Item(20,5).Value = 62

$ws.Cells.Item(21,1).Value = "HRAD"
$ws.Cells.Item(21,2).Value = 2025
$ws.Cells.Item(21,3).Value = 9
$ws.Cells.Item(21,4).Value = 126
$ws.Cells.Item(21,5).Value = 60

$ws.Cells.Item(22,1).Value = "HRJP"
$ws.Cells.Item(22,2).Value = 2025
$ws.Cells.Item(22,3).Value = 10
$ws.Cells.Item(22,4).Value = 106
$ws.Cells.Item(22,5).Value = 59

$ws.Cells.Item(23,1).Value = "HRJP"
$ws.Cells.Item(23,2).Value = 2025
$ws.Cells.Item(23,3).Value = 11
$ws.Cells.Item(23,4).Value = 115
$ws.Cells.Item(23,5).Value = 54

$ws.Cells.Item(24,1).Value = "HRJP"
$ws.Cells.Item(24,2).Value = 2025
$ws.Cells.Item(24,3).Value = 1
$ws.Cells.Item(24,4).Value = 113
$ws.Cells.Item(24,5).Value = 54

$ws.Cells.Item(25,1).Value = "HRJP"
$ws.Cells.Item(25,2).Value = 2025
$ws.Cells.Item(25,3).Value = 2
$ws.Cells.Item(25,4).Value = 110
$ws.Cells.Item(25,5).Value = 53

$ws.Cells.Item(26,1).Value = "HRJP"
$ws.Cells.Item(26,2).Value = 2025
$ws.Cells.Item(26,3).Value = 3
$ws.Cells.Item(26,4).Value = 160
$ws.Cells.Item(26,5).Value = 79

$ws.Cells.Item(27,1).Value = "HRJP"
$ws.Cells.Item(27,2).Value = 2025
$ws.Cells.Item(27,3).Value = 4
$ws.Cells.Item(27,4).Value = 135
$ws.Cells.Item(27,5).Value = 56

$ws.Cells.Item(28,1).Value = "HRJP"
$ws.Cells.Item(28,2).Value = 2025
$ws.Cells.Item(28,3).Value = 5
$ws.Cells.Item(28,4).Value = 139
$ws.Cells.Item(28,5).Value = 71

$ws.Cells.Item(29,1).Value = "HRJP"
$ws.Cells.Item(29,2).Value = 2025
$ws.Cells.Item(29,3).Value = 6
$ws.Cells.Item(29,4).Value = 127
$ws.Cells.Item(29,5).Value = 62

$ws.Cells.Item(30,1).Value = "HRJP"
$ws.Cells.Item(30,2).Value = 2025
$ws.Cells.Item(30,3).Value = 7
$ws.Cells.Item(30,4).Value = 131
$ws.Cells.Item(30,5).Value = 69

$ws.Cells.Item(31,1).Value = "HRJP"
$ws.Cells.Item(31,2).Value = 2025
$ws.Cells.Item(31,3).Value = 8
$ws.Cells.Item(31,4).Value = 139
$ws.Cells.Item(31,5).Value = 67

$ws.Cells.Item(32,1).Value = "MOV"
$ws.Cells.Item(32,2).Value = 2025
$ws.Cells.Item(32,3).Value = 9
$ws.Cells.Item(32,4).Value = 108
$ws.Cells.Item(32,5).Value = 54

$ws.Cells.Item(33,1).Value = "MOV"
$ws.Cells.Item(33,2).Value = 2025
$ws.Cells.Item(33,3).Value = 10
$ws.Cells.Item(33,4).Value = 123
$ws.Cells.Item(33,5).Value = 59

$ws.Cells.Item(34,1).Value = "MOV"
$ws.Cells.Item(34,2).Value = 2025
$ws.Cells.Item(34,3).Value = 11
$ws.Cells.Item(34,4).Value = 122
$ws.Cells.Item(34,5).Value = 65

$ws.Cells.Item(35,1).Value = "MOV"
$ws.Cells.Item(35,2).Value = 2025
$ws.Cells.Item(35,3).Value = 1
$ws.Cells.Item(35,4).Value = 278
$ws.Cells.Item(35,5).Value = 104

$ws.Cells.Item(36,1).Value = "MOV"
$ws.Cells.Item(36,2).Value = 2025
$ws.Cells.Item(36,3).Value = 2
$ws.Cells.Item(36,4).Value = 242
$ws.Cells.Item(36,5).Value = 90

$ws.Cells.Item(37,1).Value = "MOV"
$ws.Cells.Item(37,2).Value = 2025
$ws.Cells.Item(37,3).Value = 3
$ws.Cells.Item(37,4).Value = 296
$ws.Cells.Item(37,5).Value = 113

$ws.Cells.Item(38,1).Value = "MOV"
$ws.Cells.Item(38,2).Value = 2025
$ws.Cells.Item(38,3).Value = 4
$ws.Cells.Item(38,4).Value = 296
$ws.Cells.Item(38,5).Value = 113

$ws.Cells.Item(39,1).Value = "MOV"
$ws.Cells.Item(39,2).Value = 2025
$ws.Cells.Item(39,3).Value = 5
$ws.Cells.Item(39,4).Value = 296
$ws.Cells.Item(39,5).Value = 119

$ws.Cells.Item(40,1).Value = "MOV"
$ws.Cells.Item(40,2).Value = 2025
$ws.Cells.Item(40,3).Value = 6
$ws.Cells.Item(40,4).Value = 286
$ws.Cells.Item(40,5).Value = 99

$ws.Cells.Item(41,1).Value = "MOV"
$ws.Cells.Item(41,2).Value = 2025
$ws.Cells.Item(41,3).Value = 7
$ws.Cells.Item(41,4).Value = 255
$ws.Cells.Item(41,5).Value = 92

$ws.Cells.Item(42,1).Value = "MOV"
$ws.Cells.Item(42,2).Value = 2025
$ws.Cells.Item(42,3).Value = 8
$ws.Cells.Item(42,4).Value = 262
$ws.Cells.Item(42,5).Value = 97

$ws.Cells.Item(43,1).Value = "MOV"
$ws.Cells.Item(43,2).Value = 2025
$ws.Cells.Item(43,3).Value = 9
$ws.Cells.Item(43,4).Value = 248
$ws.Cells.Item(43,5).Value = 91

$ws.Cells.Item(44,1).Value = "MOV"
$ws.Cells.Item(44,2).Value = 2025
$ws.Cells.Item(44,3).Value = 10
$ws.Cells.Item(44,4).Value = 266
$ws.Cells.Item(44,5).Value = 92

$ws.Cells.Item(45,1).Value = "MOV"
$ws.Cells.Item(45,2).Value = 2025
$ws.Cells.Item(45,3).Value = 11
$ws.Cells.Item(45,4).Value = 246
$ws.Cells.Item(45,5).Value = 87

$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("A47:XFD47").Select()
